$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-14"

# Update the header label for the current-year column (I1)
$ws.Range("I1").Value = "2022 (through 06-14)"

# Update July 2022 value (row 7)
$ws.Range("I7").Value = 65

# Update Total 2022 value (row 14)
$ws.Range("I14").Value = 728
